$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 121, pushing existing rows 121-138 down to 122-139
$ws.Rows.Item(121).Insert()

# Fill the fixed (repeated) columns for the new row, matching the surrounding rows
$ws.Cells.Item(121, 1).Value = 10
$ws.Cells.Item(121, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(121, 3).Value = "La Araucanía"
$ws.Cells.Item(121, 5).Value = 9
$ws.Cells.Item(121, 6).Value = "Fruta"
$ws.Cells.Item(121, 7).Value = 100108
$ws.Cells.Item(121, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(121, 9).Value = 100108004
$ws.Cells.Item(121, 10).Value = "Papaya"
$ws.Cells.Item(121, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(121, 18).Value = "Provincia del Elquí"

# Fill the new data specific to this record
$ws.Cells.Item(121, 4).Value = 45244
$ws.Cells.Item(121, 4).Style = $ws.Cells.Item(122, 4).Style
$ws.Cells.Item(121, 4).NumberFormat = $ws.Cells.Item(122, 4).NumberFormat
$ws.Cells.Item(121, 12).Value = "Primera"
$ws.Cells.Item(121, 13).Value = 80
$ws.Cells.Item(121, 14).Value = 26000
$ws.Cells.Item(121, 15).Value = 26000
$ws.Cells.Item(121, 16).Value = 26000
$ws.Cells.Item(121, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(121, 19).Value = 2600
$ws.Cells.Item(121, 20).Value = 10
